# EEG trigger and AOI updates
# Added EEG triggers based on diff_level, updated AOI co-ords

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: column A now gets its own explicit (wider) width so the
# "images" header/filenames aren't clipped; B/C/D keep their existing
# best-fit width untouched ---
$ws.Columns.Item(1).ColumnWidth = 16.59

# --- Row heights: header + data rows now a touch shorter ---
$ws.Range("A1:D7").RowHeight = 18.75

# --- Updated AOI co-ordinates (target_x / target_y) ---
$ws.Range("C2").Value = 0.16018518518518499
$ws.Range("D2").Value = -0.43611111111111101
$ws.Range("C3").Value = 0.75833333333333297
$ws.Range("D3").Value = -0.099074074074073995
$ws.Range("C4").Value = 0.42499999999999999
$ws.Range("D4").Value = -0.20740740740740701
$ws.Range("C5").Value = -0.42777777777777698
$ws.Range("D5").Value = -0.19166666666666601
$ws.Range("C6").Value = -0.50925925925925897
$ws.Range("D6").Value = 0.141666666666666
$ws.Range("C7").Value = 0.67592592592592504
$ws.Range("D7").Value = 0.29537037037037001

# The new numeric AOI cells no longer carry the old bordered/number style -
# reset C2:D7 back to the workbook's default (General) style, keeping only
# C1:D1 headers on the right-aligned numeric style.
$ws.Range("C2:D7").Style = "Normal"

# --- Selection moves to A2 ---
$ws.Range("A2").Select()
